{"js": "// Replace the 25 \"three-digit \u00f7 one-digit\" division prompts in the practice\n// table with their new values, preserving all existing run formatting\n// (font, size, etc.) by using a search-and-replace (InsertLocation.Replace)\n// on each exact original text instead of rebuilding runs.\n\nconst replacements = [\n  [\"638\u00f79=\", \"191\u00f73=\"],\n  [\"282\u00f76=\", \"901\u00f78=\"],\n  [\"615\u00f75=\", \"319\u00f78=\"],\n  [\"123\u00f74=\", \"432\u00f78=\"],\n  [\"537\u00f78=\", \"841\u00f75=\"],\n  [\"765\u00f76=\", \"771\u00f72=\"],\n  [\"723\u00f74=\", \"771\u00f73=\"],\n  [\"821\u00f75=\", \"800\u00f79=\"],\n  [\"150\u00f75=\", \"932\u00f75=\"],\n  [\"713\u00f73=\", \"639\u00f79=\"],\n  [\"564\u00f72=\", \"640\u00f75=\"],\n  [\"943\u00f79=\", \"614\u00f77=\"],\n  [\"486\u00f74=\", \"702\u00f77=\"],\n  [\"710\u00f79=\", \"394\u00f79=\"],\n  [\"795\u00f73=\", \"112\u00f72=\"],\n  [\"689\u00f72=\", \"758\u00f73=\"],\n  [\"791\u00f72=\", \"876\u00f79=\"],\n  [\"632\u00f77=\", \"449\u00f76=\"],\n  [\"897\u00f74=\", \"340\u00f77=\"],\n  [\"839\u00f73=\", \"553\u00f78=\"],\n  [\"922\u00f77=\", \"510\u00f73=\"],\n  [\"791\u00f75=\", \"777\u00f75=\"],\n  [\"143\u00f75=\", \"805\u00f79=\"],\n  [\"335\u00f78=\", \"717\u00f76=\"],\n  [\"726\u00f72=\", \"828\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit \u00f7 one-digit\" division prompts in the practice\n# table with their new values. Word's Find/Replace keeps the existing run\n# formatting (font, size, etc.) on the matched text intact.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"638\u00f79=\", \"191\u00f73=\"),\n    @(\"282\u00f76=\", \"901\u00f78=\"),\n    @(\"615\u00f75=\", \"319\u00f78=\"),\n    @(\"123\u00f74=\", \"432\u00f78=\"),\n    @(\"537\u00f78=\", \"841\u00f75=\"),\n    @(\"765\u00f76=\", \"771\u00f72=\"),\n    @(\"723\u00f74=\", \"771\u00f73=\"),\n    @(\"821\u00f75=\", \"800\u00f79=\"),\n    @(\"150\u00f75=\", \"932\u00f75=\"),\n    @(\"713\u00f73=\", \"639\u00f79=\"),\n    @(\"564\u00f72=\", \"640\u00f75=\"),\n    @(\"943\u00f79=\", \"614\u00f77=\"),\n    @(\"486\u00f74=\", \"702\u00f77=\"),\n    @(\"710\u00f79=\", \"394\u00f79=\"),\n    @(\"795\u00f73=\", \"112\u00f72=\"),\n    @(\"689\u00f72=\", \"758\u00f73=\"),\n    @(\"791\u00f72=\", \"876\u00f79=\"),\n    @(\"632\u00f77=\", \"449\u00f76=\"),\n    @(\"897\u00f74=\", \"340\u00f77=\"),\n    @(\"839\u00f73=\", \"553\u00f78=\"),\n    @(\"922\u00f77=\", \"510\u00f73=\"),\n    @(\"791\u00f75=\", \"777\u00f75=\"),\n    @(\"143\u00f75=\", \"805\u00f79=\"),\n    @(\"335\u00f78=\", \"717\u00f76=\"),\n    @(\"726\u00f72=\", \"828\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
